$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-12: Column A = symptom_group label, B = Gas-only, C = ICD+Gas, D = ICD-only
$data = @(
    @{ Row = 2;  Label = "Diseases (patient-stated)";      B = 6.5;                C = 3.1; D = 3.5 },
    @{ Row = 3;  Label = "Injuries & adverse effects";     B = 9.699999999999999;  C = 4.4; D = 4.1 },
    @{ Row = 4;  Label = "Other";                          B = 6.7;                C = 4.8; D = 7.4 },
    @{ Row = 5;  Label = "Symptom – Circulatory";          B = 10.6;               C = 6.6; D = 7.8 },
    @{ Row = 6;  Label = "Symptom – Digestive";            B = 12.6;               C = 7;   D = 5.9 },
    @{ Row = 7;  Label = "Symptom – General";              B = 5.2;                C = 3.9; D = 6.3 },
    @{ Row = 8;  Label = "Symptom – Genitourinary";        B = 5.9;                C = 4.3; D = 3.7 },
    @{ Row = 9;  Label = "Symptom – Nervous";               B = 11.4;              C = 12;  D = 10.9 },
    @{ Row = 10; Label = "Symptom – Respiratory";          B = 26.5;               C = 51.8; D = 45.2 },
    @{ Row = 11; Label = "Symptom – Skin/Hair/Nails";      B = 2.6;                C = 1.4; D = 3.3 },
    @{ Row = 12; Label = "Uncodable/Unknown";               B = 2.3;               C = 0.8; D = 2 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.Label
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
}
